# "Generate Report for Handback" -- records the handback of the zh-cn and
# de-de translations: target (.md) + handback (.xlf) files are stamped onto
# each locale sheet along with a handback timestamp, the overall status
# flips from "Ready for handoff" to "Handed back: in sync with en-US", and
# a few columns are widened to comfortably fit the (now longer) file-name /
# status strings.

$wb = $excel.ActiveWorkbook

$mdFile  = "19175915-1708-4e2e-af2b-7cfe3d5a7cba.md"
$mdUrl   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dc6c214ed87eda4aabf5cf0dc61f5282be88ea8b/e2e/19175915-1708-4e2e-af2b-7cfe3d5a7cba.md"
$zhXlf   = "19175915-1708-4e2e-af2b-7cfe3d5a7cba.598a90aa287260df44adbfcb02060d4f088db8e5.zh-cn.xlf"
$deXlf   = "19175915-1708-4e2e-af2b-7cfe3d5a7cba.598a90aa287260df44adbfcb02060d4f088db8e5.de-de.xlf"
$status  = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: status cells for both locales + widen the two status
# columns (E, F) so the longer status text is readable.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $status
$wsOverview.Range("F2").Value = $status
$wsOverview.Range("E1:F1").ColumnWidth = 29.17

# ---------------------------------------------------------------------
# zh-cn sheet: stamp the handback info onto row 2 and widen columns.
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $status
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdFile)
$wsZh.Range("J2").Value = $zhXlf
$wsZh.Range("K2").Value = "2016-10-24 09:58:14"
$wsZh.Range("C1").ColumnWidth = 29.17
$wsZh.Range("I1:J1").ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet: stamp the handback info onto row 2 and widen columns.
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $status
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdFile)
$wsDe.Range("J2").Value = $deXlf
$wsDe.Range("K2").Value = "2016-10-24 09:58:31"
$wsDe.Range("C1").ColumnWidth = 29.17
$wsDe.Range("I1:J1").ColumnWidth = 39.17
